$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1118.1578
$ws.Range("I41").Value = 925
$ws.Range("J41").Value = 1207.3077
$ws.Range("K41").Value = 925
$ws.Range("L41").Value = 1207.3077
$ws.Range("M41").Value = -485
$ws.Range("N41").Value = -2087.3077

$ws.Range("H98").Value = 2502.1365
$ws.Range("I98").Value = 2341.2222
$ws.Range("K98").Value = 2341.2222
$ws.Range("M98").Value = -843.2222000000002

$ws.Range("H122").Value = 2502.1365
$ws.Range("I122").Value = 2341.2222
$ws.Range("K122").Value = 7023.6666
$ws.Range("M122").Value = -4573.6666

$ws.Range("H132").Value = 4096.4287
$ws.Range("I132").Value = 4281.864
$ws.Range("K132").Value = 12845.592
$ws.Range("M132").Value = -10315.592

$ws.Range("H137").Value = 78901.69500000001
$ws.Range("J137").Value = 3076
$ws.Range("L137").Value = 9228
$ws.Range("N137").Value = -14328

$ws.Range("H141").Value = 2593.4614
$ws.Range("I141").Value = 2370.625
$ws.Range("J141").Value = 2950
$ws.Range("K141").Value = 7111.875
$ws.Range("L141").Value = 8850
$ws.Range("M141").Value = -1931.875
$ws.Range("N141").Value = -19210

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4821.367
$ws.Range("I32").Value = 2324.7104
$ws.Range("J32").Value = 13446.182
$ws.Range("K32").Value = 2324.7104
$ws.Range("L32").Value = 13446.182
$ws.Range("M32").Value = -2037.7104
$ws.Range("N32").Value = -14020.182

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H101").Value = 66698.5
$ws.Range("J101").Value = 66698.5
$ws.Range("L101").Value = 66698.5
$ws.Range("N101").Value = -73188.5

$ws.Range("H105").Value = 24875
$ws.Range("J105").Value = 24875
$ws.Range("L105").Value = 24875
$ws.Range("N105").Value = -31863

$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524

$ws.Range("H122").Value = 3554.7
$ws.Range("I122").Value = 3362.5715
$ws.Range("J122").Value = 4003
$ws.Range("K122").Value = 10087.7145
$ws.Range("L122").Value = 12009
$ws.Range("M122").Value = -7637.7145
$ws.Range("N122").Value = -16909

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 6079.1904
$ws.Range("I94").Value = 2277.6667
$ws.Range("K94").Value = 2277.6667
$ws.Range("M94").Value = -1826.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3174.1177
$ws.Range("I58").Value = 2100.111
$ws.Range("K58").Value = 2100.111
$ws.Range("M58").Value = -1897.111

$ws.Range("H99").Value = 3935.4
$ws.Range("I99").Value = 3254.5
$ws.Range("J99").Value = 4713.5713
$ws.Range("K99").Value = 3254.5
$ws.Range("L99").Value = 4713.5713
$ws.Range("M99").Value = -1756.5
$ws.Range("N99").Value = -7709.5713

$ws.Range("H107").Value = 55557684
$ws.Range("I107").Value = 2101.5
$ws.Range("K107").Value = 2101.5
$ws.Range("M107").Value = -181.5

$ws.Range("H126").Value = 3935.4
$ws.Range("I126").Value = 3254.5
$ws.Range("J126").Value = 4713.5713
$ws.Range("K126").Value = 9763.5
$ws.Range("L126").Value = 14140.7139
$ws.Range("M126").Value = -7293.5
$ws.Range("N126").Value = -19080.7139

$ws.Range("H134").Value = 3422.5
$ws.Range("I134").Value = 3058.353
$ws.Range("K134").Value = 9175.059000000001
$ws.Range("M134").Value = -6640.059000000001

$ws.Range("H136").Value = 3174.1177
$ws.Range("I136").Value = 2100.111
$ws.Range("K136").Value = 6300.333
$ws.Range("M136").Value = -3750.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1181.3
$ws.Range("I129").Value = 916.2
$ws.Range("J129").Value = 1446.4
$ws.Range("K129").Value = 2748.6
$ws.Range("L129").Value = 4339.200000000001
$ws.Range("M129").Value = 2251.4
$ws.Range("N129").Value = -14339.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 23000
$ws.Range("I57").Value = 14000
$ws.Range("K57").Value = 14000
$ws.Range("M57").Value = -13180

$ws.Range("H80").Value = 4780.3
$ws.Range("I80").Value = 2949.25
$ws.Range("J80").Value = 6001
$ws.Range("K80").Value = 2949.25
$ws.Range("L80").Value = 6001
$ws.Range("M80").Value = -1951.25
$ws.Range("N80").Value = -7997

$ws.Range("H83").Value = 4780.3
$ws.Range("I83").Value = 2949.25
$ws.Range("J83").Value = 6001
$ws.Range("K83").Value = 14746.25
$ws.Range("L83").Value = 30005
$ws.Range("M83").Value = -9754.25
$ws.Range("N83").Value = -39989

$ws.Range("H113").Value = 1850
$ws.Range("I113").Value = 1850
$ws.Range("K113").Value = 1850
$ws.Range("M113").Value = 320

$ws.Range("H132").Value = 4759.3335
$ws.Range("I132").Value = 3664.1667
$ws.Range("J132").Value = 6949.6665
$ws.Range("K132").Value = 10992.5001
$ws.Range("L132").Value = 20848.9995
$ws.Range("M132").Value = -8462.500100000001
$ws.Range("N132").Value = -25908.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 129818.86
$ws.Range("I22").Value = 224408.5
$ws.Range("K22").Value = 224408.5
$ws.Range("M22").Value = -224113.5

$ws.Range("H27").Value = 129818.86
$ws.Range("I27").Value = 224408.5
$ws.Range("K27").Value = 224408.5
$ws.Range("M27").Value = -224301.5

$ws.Range("H40").Value = 5619.4
$ws.Range("I40").Value = 4398.8335
$ws.Range("K40").Value = 4398.8335
$ws.Range("M40").Value = -4262.8335

$ws.Range("H48").Value = 33333.332
$ws.Range("I48").Value = 30000
$ws.Range("K48").Value = 30000
$ws.Range("M48").Value = -29339

$ws.Range("H76").Value = 13192
$ws.Range("J76").Value = 13192
$ws.Range("L76").Value = 13192
$ws.Range("N76").Value = -13868

$ws.Range("H79").Value = 13192
$ws.Range("J79").Value = 13192
$ws.Range("L79").Value = 13192
$ws.Range("N79").Value = -15532

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1320.4688
$ws.Range("I122").Value = 990.875
$ws.Range("J122").Value = 2309.25
$ws.Range("K122").Value = 2972.625
$ws.Range("L122").Value = 6927.75
$ws.Range("M122").Value = -522.625
$ws.Range("N122").Value = -11827.75

$ws.Range("H136").Value = 2913.9048
$ws.Range("I136").Value = 1918.25
$ws.Range("K136").Value = 5754.75
$ws.Range("M136").Value = -3204.75
